$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.022105
$ws.Range("N2").Value = 0.066315
$ws.Range("O2").Value = 0.0007557226718989593
$ws.Range("P2").Value = 0.0007557226718989592
$ws.Range("Q2").Value = 4.452825511646667
$ws.Range("R2").Value = 40.07542960482
$ws.Range("S2").Value = 0.0003652785710772291
$ws.Range("T2").Value = 0.0003652785710772291
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.01371106452749117
$ws.Range("P3").Value = 0.01371106452749117
$ws.Range("Q3").Value = 80.78754308962846
$ws.Range("R3").Value = 727.0878878066561
$ws.Range("S3").Value = 0.006627243358949111
$ws.Range("T3").Value = 0.006627243358949112
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 28.82699233333333
$ws.Range("N4").Value = 86.480977
$ws.Range("O4").Value = 0.9855332128006099
$ws.Range("P4").Value = 0.9855332128006098
$ws.Range("Q4").Value = 5806.90191748064
$ws.Range("R4").Value = 52262.11725732576
$ws.Range("S4").Value = 0.4763575013786129
$ws.Range("T4").Value = 0.4763575013786129
$ws.Range("I5").Value = 0.1569674599353791
$ws.Range("J5").Value = 0.1569674599353792
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.022105
$ws.Range("N5").Value = 0.066315
$ws.Range("O5").Value = 0.0007557226718989593
$ws.Range("P5").Value = 0.0007557226718989592
$ws.Range("Q5").Value = 1.446050845956667
$ws.Range("R5").Value = 13.01445761361
$ws.Range("S5").Value = 0.0001186238682235576
$ws.Range("T5").Value = 0.0001186238682235576
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("O6").Value = 0.01371106452749117
$ws.Range("P6").Value = 0.01371106452749117
$ws.Range("S6").Value = 0.002152190971890368
$ws.Range("T6").Value = 0.002152190971890368
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("M7").Value = 28.82699233333333
$ws.Range("N7").Value = 86.480977
$ws.Range("O7").Value = 0.9855332128006099
$ws.Range("P7").Value = 0.9855332128006098
$ws.Range("Q7").Value = 1885.785869712871
$ws.Range("R7").Value = 16972.07282741584
$ws.Range("S7").Value = 0.1546966450952652
$ws.Range("T7").Value = 0.1546966450952652
$ws.Range("G8").Value = 60.43484133333334
$ws.Range("H8").Value = 181.304524
$ws.Range("I8").Value = 0.1450120099461104
$ws.Range("J8").Value = 0.1450120099461104
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.022105
$ws.Range("N8").Value = 0.066315
$ws.Range("O8").Value = 0.0007557226718989593
$ws.Range("P8").Value = 0.0007557226718989592
$ws.Range("Q8").Value = 1.335912167673333
$ws.Range("R8").Value = 12.02320950906
$ws.Range("S8").Value = 0.000109588863613913
$ws.Range("T8").Value = 0.000109588863613913
$ws.Range("G9").Value = 60.43484133333334
$ws.Range("H9").Value = 181.304524
$ws.Range("I9").Value = 0.1450120099461104
$ws.Range("J9").Value = 0.1450120099461104
$ws.Range("O9").Value = 0.01371106452749117
$ws.Range("P9").Value = 0.01371106452749117
$ws.Range("Q9").Value = 24.23743340662756
$ws.Range("R9").Value = 218.136900659648
$ws.Range("S9").Value = 0.00198826902563231
$ws.Range("T9").Value = 0.00198826902563231
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 28.82699233333333
$ws.Range("N10").Value = 86.480977
$ws.Range("O10").Value = 0.9855332128006099
$ws.Range("P10").Value = 0.9855332128006098
$ws.Range("Q10").Value = 1742.154707782217
$ws.Range("R10").Value = 15679.39237003995
$ws.Range("S10").Value = 0.1429141520568641
$ws.Range("T10").Value = 0.1429141520568641
$ws.Range("G11").Value = 89.46554166666668
$ws.Range("H11").Value = 268.396625
$ws.Range("I11").Value = 0.2146705068098712
$ws.Range("J11").Value = 0.2146705068098712
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.022105
$ws.Range("N11").Value = 0.066315
$ws.Range("O11").Value = 0.0007557226718989593
$ws.Range("P11").Value = 0.0007557226718989592
$ws.Range("Q11").Value = 1.977635798541667
$ws.Range("R11").Value = 17.798722186875
$ws.Range("S11").Value = 0.0001622313689842596
$ws.Range("T11").Value = 0.0001622313689842596
$ws.Range("G12").Value = 89.46554166666668
$ws.Range("H12").Value = 268.396625
$ws.Range("I12").Value = 0.2146705068098712
$ws.Range("J12").Value = 0.2146705068098712
$ws.Range("O12").Value = 0.01371106452749117
$ws.Range("P12").Value = 0.01371106452749117
$ws.Range("Q12").Value = 35.88021512911111
$ws.Range("R12").Value = 322.921936162
$ws.Range("S12").Value = 0.002943361171019376
$ws.Range("T12").Value = 0.002943361171019376
$ws.Range("G13").Value = 89.46554166666668
$ws.Range("H13").Value = 268.396625
$ws.Range("I13").Value = 0.2146705068098712
$ws.Range("J13").Value = 0.2146705068098712
$ws.Range("M13").Value = 28.82699233333333
$ws.Range("N13").Value = 86.480977
$ws.Range("O13").Value = 0.9855332128006099
$ws.Range("P13").Value = 0.9855332128006098
$ws.Range("Q13").Value = 2579.022483722514
$ws.Range("R13").Value = 23211.20235350263
$ws.Range("S13").Value = 0.2115649142698676
$ws.Range("T13").Value = 0.2115649142698676
